$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.168.83'
$ws.Range("E2").Value = '  -3.53%  '

$ws.Range("D3").Value = '1.971.95'
$ws.Range("E3").Value = '  -5.48%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '328.61'
$ws.Range("E5").Value = '  -3.64%  '

$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").Value = '0.5015'
$ws.Range("E7").Value = '  -4.99%  '

$ws.Range("D8").Value = '0.4231'
$ws.Range("E8").Value = '  -3.38%  '

$ws.Range("D9").Value = '52.91'
$ws.Range("E9").Value = '  -3.61%  '

$ws.Range("D10").Value = '0.09250'
$ws.Range("E10").Value = '  -0.99%  '

$ws.Range("D11").Value = '1.104'
$ws.Range("E11").Value = '  -5.90%  '

$ws.Range("D12").Value = '23.08'
$ws.Range("E12").Value = '  -5.75%  '

$ws.Range("D13").Value = '2.005.76'
$ws.Range("E13").Value = '  -1.92%  '

$ws.Range("D14").Value = '7.923'
$ws.Range("E14").Value = '  -6.65%  '

$ws.Range("D15").Value = '6.458'
$ws.Range("E15").Value = '  -5.81%  '

$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.13%  '

$ws.Range("D17").Value = '0.00001107'
$ws.Range("E17").Value = '  -4.37%  '

$ws.Range("D18").Value = '91.86'
$ws.Range("E18").Value = '  -9.50%  '

$ws.Range("D19").Value = '0.06722'
$ws.Range("E19").Value = '  +0.16%  '

$ws.Range("D20").Value = '19.40'
$ws.Range("E20").Value = '  -7.67%  '

$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  +0.31%  '

$ws.Range("D22").Value = '5.989'
$ws.Range("E22").Value = '  -4.74%  '

$ws.Range("D23").Value = '29.209.05'
$ws.Range("E23").Value = '  -3.43%  '

$ws.Range("D24").Value = '12.15'
$ws.Range("E24").Value = '  -1.97%  '

$ws.Range("D25").Value = '2.287'
$ws.Range("E25").Value = '  -1.37%  '

$ws.Range("D26").Value = '2.232.04'
$ws.Range("E26").Value = '  -3.15%  '

$ws.Range("D27").Value = '20.68'
$ws.Range("E27").Value = '  -4.97%  '

$ws.Range("D28").Value = '156.63'
$ws.Range("E28").Value = '  -3.55%  '

$ws.Range("D29").Value = '6.234'
$ws.Range("E29").Value = '  -9.26%  '

$ws.Range("D30").Value = '2.273'
$ws.Range("E30").Value = '  -8.40%  '

$ws.Range("D31").Value = '126.88'
$ws.Range("E31").Value = '  -4.93%  '

$ws.Range("D32").Value = '1.048'
$ws.Range("E32").Value = '  -7.14%  '

$ws.Range("D33").Value = '0.09873'
$ws.Range("E33").Value = '  -5.72%  '

$ws.Range("D34").Value = '1.542'
$ws.Range("E34").Value = '  -7.03%  '

$ws.Range("D35").Value = '5.805'
$ws.Range("E35").Value = '  -7.10%  '

$ws.Range("D36").Value = '3.674'
$ws.Range("E36").Value = '  -6.13%  '

$ws.Range("D37").Value = '0.02434'
$ws.Range("E37").Value = '  -6.87%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").Value = '9.055'
$ws.Range("E38").Value = '  -9.61%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '1.303'
$ws.Range("E39").Value = '  -2.75%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06372'
$ws.Range("E40").Value = '  -5.35%  '

$ws.Range("D41").Value = '0.6470'
$ws.Range("E41").Value = '  -6.74%  '

$ws.Range("D42").Value = '11.48'
$ws.Range("E42").Value = '  -8.62%  '

$ws.Range("D43").Value = '0.1993'
$ws.Range("E43").Value = '  -9.27%  '

$ws.Range("D44").Value = '1.005'
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("D45").Value = '0.6266'
$ws.Range("E45").Value = '  -6.87%  '

$ws.Range("D46").Value = '13.45'
$ws.Range("E46").Value = '  -6.13%  '

$ws.Range("D47").Value = '2.201'
$ws.Range("E47").Value = '  -7.70%  '

$ws.Range("D48").Value = '1.294'
$ws.Range("E48").Value = '  +0.58%  '

$ws.Range("D49").Value = '3.472'
$ws.Range("E49").Value = '  -4.52%  '

$ws.Range("D50").Value = '0.00000000329'
$ws.Range("E50").Value = '  -4.16%  '

$ws.Range("D51").Value = '0.06988'
$ws.Range("E51").Value = '  -3.54%  '
